$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range
$boldLabel = "Meta description"
$restOfText = ": Discover the Book of Tombs online slot game with 5 reels, 10 fixed paylines, high volatility rate, and free spins function. Play for free and read our review."
$metaRange.Text = $boldLabel + $restOfText

# Bold just the "Meta description" label, leaving the rest of the sentence
# in normal formatting (this naturally splits the paragraph into separate
# runs with the correct formatting boundaries).
$boldRange = $d.Range($metaRange.Start, $metaRange.Start + $boldLabel.Length)
$boldRange.Bold = 1

# ---------------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph that used to sit right
#    before the final (italic) meta-description paragraph at the bottom of
#    the document.
# ---------------------------------------------------------------------------
$oldTitleText = "Play Book of Tombs Online Slot for Free - Game Review"
$bottomTitlePara = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text.Contains($oldTitleText)) {
        $bottomTitlePara = $candidate
        break
    }
}
$bottomTitlePara.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph - the old meta
#    description copy - with the new image-prompt text, keeping its
#    existing italic run formatting intact.
# ---------------------------------------------------------------------------
$oldDescription = "Discover the Book of Tombs online slot game with 5 reels, 10 fixed paylines, high volatility rate, and free spins function. Play for free and read our review."
$newDescription = 'Create a cartoon-style feature image for the game "Book of Tombs" that features a happy Maya warrior with glasses. The image should be vibrant and eye-catching, using warm colors to evoke the Egyptian theme of the game. The Maya warrior should be shown holding the book of the Pharaoh, with a confident expression on their face. The background should feature the pyramids and other Egyptian landmarks, with the logo of the game prominently displayed. It should convey a sense of adventure and excitement, making players eager to dive into the game and uncover the treasures that await them.'

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$idx = $lastRange.Text.IndexOf($oldDescription)
$target = $d.Range($lastRange.Start + $idx, $lastRange.Start + $idx + $oldDescription.Length)
$target.Text = $newDescription
